# Cck-Cckbr LR-pair sheet refresh ("update scripts wuth new tpm"):
# the underlying TPM recompute changed which sending/target cluster
# combinations are reported. The table is now the cross join of all 4
# sending clusters (ECs, FAPs, MuSCs, Resolving-Mac) against just the
# 2 target clusters (ECs, FAPs) -- 8 data rows total (was 4 rows, all
# targeting FAPs only). Column order matches the header row:
# A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T the various expression / specificity metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  ,@("ECs","Cck","Cckbr","ECs",3,1,0.5218183333333334,1.565455,0.11537035205514,0.11537035205514,1,0.3333333333333333,0.09168033333333332,0.275041,0.1820807603088688,0.1820807603088689,0.04784047873944444,0.430564308655,0.02100672141930176,0.02100672141930176)
  ,@("ECs","Cck","Cckbr","FAPs",3,1,0.5218183333333334,1.565455,0.11537035205514,0.11537035205514,3,1,0.4118343333333334,1.235503,0.8179192396911312,0.8179192396911312,0.2149027054294445,1.934124348865,0.09436363063583826,0.09436363063583826)
  ,@("FAPs","Cck","Cckbr","ECs",3,1,1.300224333333333,3.900673,0.2874704269761693,0.2874704269761693,1,0.3333333333333333,0.09168033333333332,0.275041,0.1820807603088688,0.1820807603088689,0.1192050002881111,1.072845002593,0.05234283391013606,0.05234283391013607)
  ,@("FAPs","Cck","Cckbr","FAPs",3,1,1.300224333333333,3.900673,0.2874704269761693,0.2874704269761693,3,1,0.4118343333333334,1.235503,0.8179192396911312,0.8179192396911312,0.5354770215021111,4.819293193519,0.2351275930660333,0.2351275930660333)
  ,@("MuSCs","Cck","Cckbr","ECs",3,1,2.600134333333333,7.800402999999999,0.5748713570699702,0.5748713570699702,1,0.3333333333333333,0.09168033333333332,0.275041,0.1820807603088688,0.1820807603088689,0.2383811823914444,2.145430641523,0.1046730137750914,0.1046730137750914)
  ,@("MuSCs","Cck","Cckbr","FAPs",3,1,2.600134333333333,7.800402999999999,0.5748713570699702,0.5748713570699702,3,1,0.4118343333333334,1.235503,0.8179192396911312,0.8179192396911312,1.070824589745444,9.637421307708999,0.4701983432948788,0.4701983432948788)
  ,@("Resolving-Mac","Cck","Cckbr","ECs",1,0.3333333333333333,0.1008076666666667,0.302423,0.02228786389872057,0.02228786389872057,1,0.3333333333333333,0.09168033333333332,0.275041,0.1820807603088688,0.1820807603088689,0.009242080482555554,0.083178724343,0.00405819120433963,0.004058191204339631)
  ,@("Resolving-Mac","Cck","Cckbr","FAPs",1,0.3333333333333333,0.1008076666666667,0.302423,0.02228786389872057,0.02228786389872057,3,1,0.4118343333333334,1.235503,0.8179192396911312,0.8179192396911312,0.04151605819655556,0.373644523769,0.01822967269438094,0.01822967269438094)
)

$r = 2
foreach ($row in $rows) {
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws.Cells.Item($r, $c + 1).Value = $row[$c]
  }
  $r = $r + 1
}
